# Adds a "Custom Identifier" column at the front of the sheet (new column A),
# shifting the existing Street/Number/Neighborhood/Postal Code/State/City
# table one column to the right (now B:G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at A; existing data (A:F) shifts to B:G.
$ws.Columns("A:A").Insert()

# Bring over the header/data cell formatting from the (now shifted) first
# data column so the new column visually matches the rest of the table
# (bold header style, plain data style) instead of being left unformatted.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# New header + placeholder content.
$ws.Range("A1").Value = "Custom Identifier"
$ws.Range("A2").Value = "Enter with your data from here..."

# Match the author's final selection/active cell.
[void]$ws.Range("A2").Select()

# Page setup (paper size / orientation) as configured in the committed copy.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
